$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-20 18:18:28"
$ws.Range("E3").Value = "2026-02-20 18:18:31"
$ws.Range("E4").Value = "2026-02-20 18:18:33"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "54%"
$ws.Range("J4").Value = "1022.1 hPa"
$ws.Range("O4").Value = "10.6 °C"
$ws.Range("E5").Value = "2026-02-20 18:18:35"
$ws.Range("O5").Value = "-4.5 °C"
$ws.Range("E6").Value = "2026-02-20 18:18:38"
$ws.Range("J6").Value = "1022.1 hPa"
$ws.Range("E7").Value = "2026-02-20 18:18:40"
$ws.Range("J7").Value = "1022.0 hPa"
$ws.Range("E8").Value = "2026-02-20 18:18:43"
$ws.Range("E9").Value = "2026-02-20 18:18:45"
$ws.Range("E10").Value = "2026-02-20 18:18:47"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = "77%"
$ws.Range("E11").Value = "2026-02-20 18:18:50"
$ws.Range("E12").Value = "2026-02-20 18:18:52"
$ws.Range("E13").Value = "2026-02-20 18:18:54"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "42%"
$ws.Range("J13").Value = "1022.9 hPa"
$ws.Range("E14").Value = "2026-02-20 18:18:57"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "53%"
$ws.Range("O14").Value = "12.4 °C"
$ws.Range("E15").Value = "2026-02-20 18:18:59"
$ws.Range("E16").Value = "2026-02-20 18:19:02"
$ws.Range("M16").Value = "0.2 °C 17:59 TU"
$ws.Range("O16").Value = "-3.6 °C"
$ws.Range("E17").Value = "2026-02-20 18:19:04"
$ws.Range("E18").Value = "2026-02-20 18:19:07"
$ws.Range("H18").NumberFormat = "@"
$ws.Range("H18").Value = "75%"
$ws.Range("J18").Value = "1022.4 hPa"
$ws.Range("E19").Value = "2026-02-20 18:19:09"
$ws.Range("E20").Value = "2026-02-20 18:19:11"
$ws.Range("E21").Value = "2026-02-20 18:19:14"
$ws.Range("J21").Value = "1022.0 hPa"
$ws.Range("E22").Value = "2026-02-20 18:19:16"
$ws.Range("O22").Value = "-4.3 °C"
$ws.Range("E23").Value = "2026-02-20 18:19:19"
$ws.Range("O23").Value = "-5.3 °C"
$ws.Range("E24").Value = "2026-02-20 18:19:21"
$ws.Range("O24").Value = "9.6 °C"
$ws.Range("E25").Value = "2026-02-20 18:19:24"
$ws.Range("E26").Value = "2026-02-20 18:19:26"
$ws.Range("J26").Value = "1021.1 hPa"
$ws.Range("E27").Value = "2026-02-20 18:19:29"
$ws.Range("E28").Value = "2026-02-20 18:19:31"
$ws.Range("J28").Value = "1022.4 hPa"
$ws.Range("E29").Value = "2026-02-20 18:19:34"
$ws.Range("H29").NumberFormat = "@"
$ws.Range("H29").Value = "72%"
$ws.Range("L29").Value = "25.6 km/h - 0º 17:42 TU"
$ws.Range("O29").Value = "9.4 °C"
$ws.Range("E30").Value = "2026-02-20 18:19:36"
$ws.Range("J30").Value = "1021.7 hPa"
$ws.Range("E31").Value = "2026-02-20 18:19:38"
$ws.Range("J31").Value = "1020.9 hPa"
$ws.Range("E32").Value = "2026-02-20 18:19:40"
$ws.Range("H32").NumberFormat = "@"
$ws.Range("H32").Value = "81%"
$ws.Range("E33").Value = "2026-02-20 18:19:43"
$ws.Range("J33").Value = "1022.3 hPa"
$ws.Range("E34").Value = "2026-02-20 18:19:45"
$ws.Range("H34").NumberFormat = "@"
$ws.Range("H34").Value = "45%"
$ws.Range("O34").Value = "0.4 °C"
$ws.Range("E35").Value = "2026-02-20 18:19:48"
$ws.Range("J35").Value = "1026.4 hPa"
$ws.Range("E36").Value = "2026-02-20 18:19:50"
$ws.Range("J36").Value = "1022.0 hPa"
$ws.Range("E37").Value = "2026-02-20 18:19:52"
$ws.Range("J37").Value = "1023.8 hPa"
$ws.Range("E38").Value = "2026-02-20 18:19:55"
$ws.Range("H38").NumberFormat = "@"
$ws.Range("H38").Value = "66%"
$ws.Range("E39").Value = "2026-02-20 18:19:57"
$ws.Range("H39").NumberFormat = "@"
$ws.Range("H39").Value = "48%"
$ws.Range("O39").Value = "-2.9 °C"
$ws.Range("E40").Value = "2026-02-20 18:20:00"
$ws.Range("H40").NumberFormat = "@"
$ws.Range("H40").Value = "35%"
$ws.Range("J40").Value = "1022.9 hPa"
$ws.Range("E41").Value = "2026-02-20 18:20:02"
$ws.Range("J41").Value = "1022.7 hPa"
$ws.Range("O41").Value = "13.4 °C"
$ws.Range("E42").Value = "2026-02-20 18:20:05"
$ws.Range("H42").NumberFormat = "@"
$ws.Range("H42").Value = "67%"
$ws.Range("O42").Value = "10.0 °C"
$ws.Range("E43").Value = "2026-02-20 18:20:07"
$ws.Range("O43").Value = "5.0 °C"
$ws.Range("E44").Value = "2026-02-20 18:20:10"
$ws.Range("E45").Value = "2026-02-20 18:20:12"
$ws.Range("O45").Value = "3.8 °C"
$ws.Range("E46").Value = "2026-02-20 18:20:15"
